$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1539.1765
$ws.Range("I19").Value = 1534.6
$ws.Range("K19").Value = 1534.6
$ws.Range("M19").Value = -1359.6
$ws.Range("H62").Value = 3673.7222
$ws.Range("I62").Value = 3478.0588
$ws.Range("K62").Value = 3478.0588
$ws.Range("M62").Value = -2854.0588
$ws.Range("H65").Value = 3673.7222
$ws.Range("I65").Value = 3478.0588
$ws.Range("K65").Value = 17390.294
$ws.Range("M65").Value = -14270.294
$ws.Range("H92").Value = 712.5
$ws.Range("I92").Value = 712.5
$ws.Range("K92").Value = 712.5
$ws.Range("M92").Value = 535.5
$ws.Range("H100").Value = 2987.375
$ws.Range("J100").Value = 3033.3333
$ws.Range("L100").Value = 3033.3333
$ws.Range("N100").Value = -4115.3333
$ws.Range("H103").Value = 1415.5555
$ws.Range("I103").Value = 1415.5555
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 4246.666499999999
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -3660.666499999999
$ws.Range("N103").ClearContents()
$ws.Range("H106").Value = 1459.2
$ws.Range("I106").Value = 1459.2
$ws.Range("K106").Value = 1459.2
$ws.Range("M106").Value = -828.2
$ws.Range("H111").Value = 690
$ws.Range("I111").Value = 762.5
$ws.Range("J111").Value = 617.5
$ws.Range("K111").Value = 2287.5
$ws.Range("L111").Value = 1852.5
$ws.Range("M111").Value = 779.5
$ws.Range("N111").Value = -7986.5
$ws.Range("H116").Value = 5267.3335
$ws.Range("I116").Value = 3333
$ws.Range("J116").Value = 6234.5
$ws.Range("K116").Value = 3333
$ws.Range("L116").Value = 6234.5
$ws.Range("M116").Value = 109
$ws.Range("N116").Value = -13118.5
$ws.Range("H125").Value = 2710.6428
$ws.Range("I125").Value = 2117.375
$ws.Range("J125").Value = 3501.6667
$ws.Range("K125").Value = 19056.375
$ws.Range("L125").Value = 31515.0003
$ws.Range("M125").Value = -16596.375
$ws.Range("N125").Value = -36435.0003
$ws.Range("H138").Value = 2717
$ws.Range("I138").Value = 1249.75
$ws.Range("J138").Value = 3010.45
$ws.Range("K138").Value = 3749.25
$ws.Range("L138").Value = 9031.349999999999
$ws.Range("M138").Value = 1390.75
$ws.Range("N138").Value = -19311.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 20500
$ws.Range("I33").Value = 20500
$ws.Range("K33").Value = 20500
$ws.Range("M33").Value = -20171
$ws.Range("H36").Value = 5000
$ws.Range("I36").Value = 5000
$ws.Range("K36").Value = 5000
$ws.Range("M36").Value = -4654
$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -35976
$ws.Range("H55").Value = 27498.125
$ws.Range("J55").Value = 27498.125
$ws.Range("L55").Value = 27498.125
$ws.Range("N55").Value = -28128.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2822.4167
$ws.Range("I99").Value = 2711.9
$ws.Range("K99").Value = 2711.9
$ws.Range("M99").Value = -1213.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1010.5
$ws.Range("I19").Value = 818.5
$ws.Range("J19").Value = 1202.5
$ws.Range("K19").Value = 818.5
$ws.Range("L19").Value = 1202.5
$ws.Range("M19").Value = -648.5
$ws.Range("N19").Value = -1542.5
$ws.Range("H22").Value = 199.66667
$ws.Range("I22").Value = 49.5
$ws.Range("K22").Value = 49.5
$ws.Range("M22").Value = 300.5
$ws.Range("H24").Value = 1010.5
$ws.Range("I24").Value = 818.5
$ws.Range("J24").Value = 1202.5
$ws.Range("K24").Value = 818.5
$ws.Range("L24").Value = 1202.5
$ws.Range("M24").Value = -648.5
$ws.Range("N24").Value = -1542.5
$ws.Range("H32").Value = 4999
$ws.Range("I32").Value = 4999
$ws.Range("K32").Value = 4999
$ws.Range("M32").Value = -4683
$ws.Range("H33").Value = 2487.5
$ws.Range("I33").Value = 2487.5
$ws.Range("K33").Value = 2487.5
$ws.Range("M33").Value = -2108.5
$ws.Range("H35").Value = 601441.8
$ws.Range("I35").Value = 601441.8
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 601441.8
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -601147.8
$ws.Range("N35").ClearContents()
$ws.Range("H41").Value = 8039.3335
$ws.Range("J41").Value = 24000
$ws.Range("L41").Value = 24000
$ws.Range("N41").Value = -24856
$ws.Range("H86").Value = 9043.700000000001
$ws.Range("I86").Value = 8940.166999999999
$ws.Range("J86").Value = 9199
$ws.Range("K86").Value = 8940.166999999999
$ws.Range("L86").Value = 9199
$ws.Range("M86").Value = -7817.166999999999
$ws.Range("N86").Value = -11445
$ws.Range("H89").Value = 9043.700000000001
$ws.Range("I89").Value = 8940.166999999999
$ws.Range("J89").Value = 9199
$ws.Range("K89").Value = 44700.835
$ws.Range("L89").Value = 45995
$ws.Range("M89").Value = -39084.835
$ws.Range("N89").Value = -57227

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1000
$ws.Range("J59").Value = 1000
$ws.Range("L59").Value = 3000
$ws.Range("N59").Value = -4080
$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 6000
$ws.Range("M68").Value = -5189
$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 18000
$ws.Range("M71").Value = -13944
$ws.Range("H117").Value = 447.625
$ws.Range("I117").Value = 472.66666
$ws.Range("J117").Value = 432.6
$ws.Range("K117").Value = 1417.99998
$ws.Range("L117").Value = 1297.8
$ws.Range("M117").Value = 2024.00002
$ws.Range("N117").Value = -8181.8
$ws.Range("H122").Value = 2165
$ws.Range("J122").Value = 2165
$ws.Range("L122").Value = 19485
$ws.Range("N122").Value = -24385
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 430714
$ws.Range("I7").Value = 3000000
$ws.Range("J7").Value = 2499.6667
$ws.Range("K7").Value = 3000000
$ws.Range("L7").Value = 2499.6667
$ws.Range("M7").Value = -2999888
$ws.Range("N7").Value = -2723.6667
$ws.Range("H8").Value = 430714
$ws.Range("I8").Value = 3000000
$ws.Range("J8").Value = 2499.6667
$ws.Range("K8").Value = 3000000
$ws.Range("L8").Value = 2499.6667
$ws.Range("M8").Value = -2999861
$ws.Range("N8").Value = -2777.6667
$ws.Range("H43").Value = 11841.667
$ws.Range("I43").Value = 5025
$ws.Range("K43").Value = 5025
$ws.Range("M43").Value = -4874
$ws.Range("H57").Value = 22994
$ws.Range("I57").Value = 12000
$ws.Range("K57").Value = 12000
$ws.Range("M57").Value = -11180
$ws.Range("H102").Value = 3215
$ws.Range("I102").Value = 1402.1538
$ws.Range("K102").Value = 1402.1538
$ws.Range("M102").Value = 219.8462
$ws.Range("H122").Value = 3293.1
$ws.Range("I122").Value = 3138.1333
$ws.Range("K122").Value = 9414.3999
$ws.Range("M122").Value = -6964.3999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 7928.5713
$ws.Range("J20").Value = 7928.5713
$ws.Range("L20").Value = 7928.5713
$ws.Range("N20").Value = -8380.5713
$ws.Range("H23").Value = 6500
$ws.Range("I23").Value = 6500
$ws.Range("K23").Value = 6500
$ws.Range("M23").Value = -6270
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H68").Value = 18059.084
$ws.Range("I68").Value = 1538.875
$ws.Range("J68").Value = 51099.5
$ws.Range("K68").Value = 1538.875
$ws.Range("L68").Value = 51099.5
$ws.Range("M68").Value = -789.875
$ws.Range("N68").Value = -52597.5
$ws.Range("H71").Value = 18059.084
$ws.Range("I71").Value = 1538.875
$ws.Range("J71").Value = 51099.5
$ws.Range("K71").Value = 7694.375
$ws.Range("L71").Value = 255497.5
$ws.Range("M71").Value = -3950.375
$ws.Range("N71").Value = -262985.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 350
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 350
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 350
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -808

Write-Output "Applied Marilith profit updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
